$d = $word.ActiveDocument

# Locate the sentence containing the typo "Giver" (should read "Given").
$rFull = $d.Content
$found = $rFull.Find.Execute("Giver the above tables solve the following queries.",
                              $false, $false, $false, $false, $false,
                              $true, 1, $false, "", 0)

$start = $rFull.Start

# Fix the typo in place: "Giver" -> "Given" (replace the "r" with "n").
$rChar = $d.Range($start + 4, $start + 5)
$rChar.Text = "n"

# The corrected word ends up split across separate runs ("Give" | "n" | rest)
# even though the three runs share identical formatting. Reproduce that by
# inserting and immediately removing temporary bookmarks at the split points,
# which forces the paragraph's runs to be cut at those offsets without
# altering the text or formatting.
$p1 = $d.Range($start + 4, $start + 4)
$d.Bookmarks.Add("tmpSplit1", $p1) | Out-Null
$d.Bookmarks("tmpSplit1").Delete()

$p2 = $d.Range($start + 5, $start + 5)
$d.Bookmarks.Add("tmpSplit2", $p2) | Out-Null
$d.Bookmarks("tmpSplit2").Delete()
